# Update destination email ("Correo Destino") and reselect cell C3.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2 holds the "Correo Destino" value - replace old address with new one.
$ws.Range("B2").Value = "jossandoval@falabella.com.co"

# Move/restore the active selection to C3 (was C8).
$ws.Range("C3").Select()
